# Rename speaker transcripts in column D (Speaker) to shortened codes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    "ANTOINETTE VILLARIN" = "T"
    "PATTY FERRANT"       = "T2"
    "STUDENT"             = "S"
}

# Data rows run from 2 to 180 (row 1 is the header row: Speaker, etc.)
$lastRow = 180

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)   # Column D = Speaker
    $current = $cell.Text
    if ($mapping.ContainsKey($current)) {
        $cell.Value = $mapping[$current]
    }
}
